$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 18 ("@jobs_internships_group") progress update:
#  - current_phase (D18): 1 -> 2
#  - last_action_date (E18): updated to new timestamp
#  - replies_count (I18): 0 -> 1
#  - replied_message_ids (M18): [] -> [14706]

$ws.Range("D18").Value = 2
$ws.Range("E18").Value = "2026-02-19T09:54:42.847070+00:00"
$ws.Range("I18").Value = 1
$ws.Range("M18").Value = "[14706]"
